# Atualiza o conteúdo das planilhas "max-arrecad" e "tx-sucesso" para
# refletir a nova ordenação de categorias empatadas, gerada por módulo
# de análise.

$wb = $excel.ActiveWorkbook

# --- Planilha "max-arrecad" ---
$wsMax = $wb.Worksheets.Item("max-arrecad")

$wsMax.Range("A2").Value = "fiq"
$wsMax.Range("A3").Value = "ficcao_cientifica"
$wsMax.Range("A4").Value = "questoes_genero"
$wsMax.Range("A5").Value = "fantasia"
$wsMax.Range("A6").Value = "humor"
$wsMax.Range("A7").Value = "folclore"
$wsMax.Range("A8").Value = "religiosidade"
$wsMax.Range("A9").Value = "terror"

$wsMax.Range("A11").Value = "jogos"
$wsMax.Range("A12").Value = "webformatos"

$wsMax.Range("A13").Value = "angelo_agostini"
$wsMax.Range("A14").Value = "hqmix"

$wsMax.Range("A16").Value = "erotismo"
$wsMax.Range("A17").Value = "zine"

# --- Planilha "tx-sucesso" ---
$wsTx = $wb.Worksheets.Item("tx-sucesso")

$wsTx.Range("A2").Value = "saloes_humor"
$wsTx.Range("A5").Value = "angelo_agostini"

$wsTx.Range("A8").Value = "erotismo"
$wsTx.Range("A9").Value = "questoes_genero"

$wsTx.Range("A18").Value = "religiosidade"
$wsTx.Range("A19").Value = "herois"
